$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.971.35'
$ws.Range("E2").Value = '  -0.28%  '
$ws.Range("D3").Value = '2.409.91'
$ws.Range("E3").Value = '  -0.38%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '562.35'
$ws.Range("E5").Value = '  +0.99%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.20'
$ws.Range("E6").Value = '  -0.95%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.529'
$ws.Range("E8").Value = '  -0.76%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.109'
$ws.Range("E9").Value = '  -0.48%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.153'
$ws.Range("E10").Value = '  -1.96%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.21'
$ws.Range("E11").Value = '  -3.47%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.349'
$ws.Range("E12").Value = '  -1.15%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '25.57'
$ws.Range("E13").Value = '  -2.65%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000173'
$ws.Range("E14").Value = '  -1.99%  '
$ws.Range("D15").Value = '2.845.70'
$ws.Range("E15").Value = '  -0.36%  '
$ws.Range("D16").Value = '61.869.94'
$ws.Range("E16").Value = '  -0.22%  '
$ws.Range("D17").Value = '2.410.02'
$ws.Range("E17").Value = '  -0.38%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.22'
$ws.Range("E18").Value = '  +0.97%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '321.90'
$ws.Range("E19").Value = '  -0.95%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.13'
$ws.Range("E20").Value = '  -1.80%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.81'
$ws.Range("E21").Value = '  +1.07%  '
$ws.Range("E22").Value = '  -0.04%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '66.30'
$ws.Range("E23").Value = '  +1.79%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.73'
$ws.Range("E24").Value = '  -1.53%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.64'
$ws.Range("E25").Value = '  -4.28%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '559.86'
$ws.Range("E26").Value = '  -3.63%  '
$ws.Range("D27").Value = '2.532.78'
$ws.Range("E27").Value = '  -0.12%  '
$ws.Range("E28").Value = '  +0.35%  '
$ws.Range("D29").Value = '0.0₃0928'
$ws.Range("E29").Value = '  -1.79%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.17'
$ws.Range("E30").Value = '  -2.15%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.39'
$ws.Range("E31").Value = '  -5.40%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.147'
$ws.Range("E32").Value = '  -0.84%  '
$ws.Range("E33").Value = '  +0.02%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.51'
$ws.Range("E34").Value = '  -3.87%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  +0.06%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.75'
$ws.Range("E36").Value = '  -1.49%  '
$ws.Range("B37").Value = 'PolygonEcosystemToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.379'
$ws.Range("E37").Value = '  -1.53%  '
$ws.Range("B38").Value = 'Monero'
$ws.Range("C38").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '152.68'
$ws.Range("E38").Value = '  +2.46%  '
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.43'
$ws.Range("E39").Value = '  -4.90%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.52'
$ws.Range("E40").Value = '  -1.40%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.82'
$ws.Range("E41").Value = '  -2.43%  '
$ws.Range("E42").Value = '  -0.12%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '147.44'
$ws.Range("E43").Value = '  -2.77%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.23'
$ws.Range("E44").Value = '  -3.95%  '
$ws.Range("E45").Value = '  -1.13%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0528'
$ws.Range("E46").Value = '  -3.20%  '
$ws.Range("B47").Value = 'Mantle'
$ws.Range("C47").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.592'
$ws.Range("E47").Value = '  +0.42%  '
$ws.Range("B48").Value = 'InjectiveProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '19.80'
$ws.Range("E48").Value = '  -2.97%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0920'
$ws.Range("E49").Value = '  +0.02%  '
$ws.Range("E50").Value = '  -1.27%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '11.54'
$ws.Range("E51").Value = '  +0.48%  '
